$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their text formatting so values
# like "1.00" or "5.41" are not reinterpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "54.848.21"
$ws.Range("D3").Value = "2.347.88"
$ws.Range("D5").Value = "472.46"
$ws.Range("D6").Value = "143.20"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "0.510"
$ws.Range("D9").Value = "2.346.45"
$ws.Range("D10").Value = "0.0959"
$ws.Range("D11").Value = "5.41"
$ws.Range("D14").Value = "2.756.59"
$ws.Range("D15").Value = "54.883.56"
$ws.Range("D16").Value = "19.95"
$ws.Range("D18").Value = "2.349.84"
$ws.Range("D19").Value = "4.53"
$ws.Range("D20").Value = "311.96"
$ws.Range("D22").Value = "1.00"
$ws.Range("D24").Value = "55.96"
$ws.Range("D28").Value = "2.448.28"
$ws.Range("D31").Value = "0.0₃0748"
$ws.Range("D32").Value = "147.03"
$ws.Range("D33").Value = "17.98"
$ws.Range("D35").Value = "5.01"
$ws.Range("D37").Value = "3.55"
$ws.Range("D38").Value = "0.815"
$ws.Range("D41").Value = "1.34"
$ws.Range("D42").Value = "3.34"
$ws.Range("D43").Value = "0.0951"
$ws.Range("D47").Value = "250.74"
$ws.Range("D48").Value = "0.0220"
$ws.Range("D49").Value = "4.38"
$ws.Range("D50").Value = "16.64"
$ws.Range("D51").Value = "1.777.82"

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("E6").Value = "  -2.21%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -5.97%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  -6.48%  "
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("E18").Value = "  -5.69%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("E24").Value = "  -4.43%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("E28").Value = "  -5.59%  "
$ws.Range("E29").Value = "  -6.31%  "
$ws.Range("E31").Value = "  -5.07%  "
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("E37").Value = "  -4.49%  "
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  -4.50%  "
$ws.Range("E43").Value = "  +3.23%  "
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("E45").Value = "  -6.59%  "
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("E49").Value = "  -8.13%  "
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("E51").Value = "  -5.92%  "

# --- Rows 39/40 swapped (OKB and FirstDigitalUSD changed rank order) ---
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "33.32"
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "0.997"
$ws.Range("E40").Value = "  +0.20%  "

